$d = $word.ActiveDocument

# 1. Insert the new "Netlify" backstory sentence right after
#    "...will be handled by Heroku. " and before the second "Heroku" run.
#    (kept inside the uncoloured run so the inserted text has no w:color)
$d.Content.Find.Execute(
    "will be handled by Heroku. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "will be handled by Heroku. Our initial plan was to have our web hosting with Netlify. However, Netlify only offers a service for static web pages which would not be appropriate for our application. We require the hosting of a dynamic web application, and after researching this we found that Heroku is the best option for this. ",
    2
) | Out-Null

# 2. Rewrite the description of Heroku so it compares it to Netlify.
$d.Content.Find.Execute(
    "a cloud-based application hosting service which allows developers to run their",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "a cloud-based application hosting service similar to Netlify, but it allows developers to run their",
    2
) | Out-Null
